# Update countries & provincias Spain
# This script refreshes the COVID-19 country statistics table: it updates the
# "last updated" timestamp, refreshes numeric stats for a number of countries
# (some of which overtook others in rank, shifting their row position), while
# keeping the rest of the sheet untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 22 de Octubre de 2020 a las 11:42"

# --- Row 4: Estados Unidos (rank unchanged, fresh numbers) ---
$ws.Range("B4").Value = 8585748
$ws.Range("C4").Value = 929
$ws.Range("D4").Value = 5603025
$ws.Range("E4").Value = 2755304
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 10
$ws.Range("H4").Value = 227419

# --- Row 20: Banglades (rank unchanged, fresh numbers) ---
$ws.Range("B20").Value = 394827
$ws.Range("C20").Value = 1696
$ws.Range("D20").Value = 310532
$ws.Range("E20").Value = 78548
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 24
$ws.Range("H20").Value = 5747

# --- Row 22: Indonesia (rank unchanged, fresh numbers) ---
$ws.Range("B22").Value = 377541
$ws.Range("C22").Value = 4432
$ws.Range("D22").Value = 301006
$ws.Range("E22").Value = 63576
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 102
$ws.Range("H22").Value = 12959

# --- Rows 31-33: Polonia overtakes Chequia and Canada ---
# Row 31 now shows Polonia with fresh numbers
$ws.Range("A31").Value = "Polonia"
$ws.Range("B31").Value = 214686
$ws.Range("C31").Value = 12107
$ws.Range("D31").Value = 102204
$ws.Range("E31").Value = 108463
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 168
$ws.Range("H31").Value = 4019

# Row 32 now shows Chequia (its previous, unchanged numbers)
$ws.Range("A32").Value = "Chequia"
$ws.Range("B32").Value = 208915
$ws.Range("C32").Value = 0
$ws.Range("D32").Value = 83136
$ws.Range("E32").Value = 124040
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 1739

# Row 33 now shows Canada (its previous, unchanged numbers)
$ws.Range("A33").Value = "Canada"
$ws.Range("B33").Value = 206360
$ws.Range("C33").Value = 406
$ws.Range("D33").Value = 173748
$ws.Range("E33").Value = 22783
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = 3
$ws.Range("H33").Value = 9829

# --- Row 86: Croacia (rank unchanged, fresh numbers) ---
$ws.Range("B86").Value = 29850
$ws.Range("C86").Value = 1563
$ws.Range("D86").Value = 22064
$ws.Range("E86").Value = 7380
$ws.Range("F86").Value = 0
$ws.Range("G86").Value = 13
$ws.Range("H86").Value = 406

# --- Row 102: Finlandia (rank unchanged, fresh numbers) ---
$ws.Range("B102").Value = 14255
$ws.Range("C102").Value = 184
$ws.Range("D102").Value = 9800
$ws.Range("E102").Value = 4100
$ws.Range("F102").Value = 0
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 355

# --- Row 106: Guinea (rank unchanged, fresh numbers) ---
$ws.Range("B106").Value = 11635
$ws.Range("C106").Value = 36
$ws.Range("D106").Value = 10474
$ws.Range("E106").Value = 1090
$ws.Range("F106").Value = 0
$ws.Range("G106").Value = 1
$ws.Range("H106").Value = 71

# --- Row 109: Consejo Danes para los Refugiados (rank unchanged, fresh numbers) ---
$ws.Range("B109").Value = 11078
$ws.Range("C109").Value = 12
$ws.Range("D109").Value = 10362
$ws.Range("E109").Value = 413
$ws.Range("F109").Value = 0
$ws.Range("G109").Value = 0
$ws.Range("H109").Value = 303

# --- Rows 115-117: Lituania overtakes Jamaica and Angola ---
# Row 115 now shows Lituania with fresh numbers
$ws.Range("A115").Value = "Lituania"
$ws.Range("B115").Value = 8663
$ws.Range("C115").Value = 424
$ws.Range("D115").Value = 3773
$ws.Range("E115").Value = 4765
$ws.Range("F115").Value = 0
$ws.Range("G115").Value = 5
$ws.Range("H115").Value = 125

# Row 116 now shows Jamaica (its previous, unchanged numbers)
$ws.Range("A116").Value = "Jamaica"
$ws.Range("B116").Value = 8445
$ws.Range("C116").Value = 0
$ws.Range("D116").Value = 4016
$ws.Range("E116").Value = 4255
$ws.Range("F116").Value = 0
$ws.Range("G116").Value = 0
$ws.Range("H116").Value = 174

# Row 117 now shows Angola (its previous, unchanged numbers)
$ws.Range("A117").Value = "Angola"
$ws.Range("B117").Value = 8338
$ws.Range("C117").Value = 0
$ws.Range("D117").Value = 3040
$ws.Range("E117").Value = 5043
$ws.Range("F117").Value = 0
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 255

# --- Row 124: Sri Lanka (rank unchanged, only D/E refreshed) ---
$ws.Range("D124").Value = 3561
$ws.Range("E124").Value = 2404

# --- Rows 145-148: Letonia overtakes Somalia, Guyana and Principado de Andorra ---
# Row 145 now shows Letonia with fresh numbers
$ws.Range("A145").Value = "Letonia"
$ws.Range("B145").Value = 3958
$ws.Range("C145").Value = 161
$ws.Range("D145").Value = 1357
$ws.Range("E145").Value = 2552
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 2
$ws.Range("H145").Value = 49

# Row 146 now shows Somalia (its previous, unchanged numbers)
$ws.Range("A146").Value = "Somalia"
$ws.Range("B146").Value = 3897
$ws.Range("C146").Value = 7
$ws.Range("D146").Value = 3166
$ws.Range("E146").Value = 629
$ws.Range("F146").Value = 0
$ws.Range("G146").Value = 1
$ws.Range("H146").Value = 102

# Row 147 now shows Guyana (its previous, unchanged numbers)
$ws.Range("A147").Value = "Guyana"
$ws.Range("B147").Value = 3850
$ws.Range("C147").Value = 0
$ws.Range("D147").Value = 2839
$ws.Range("E147").Value = 895
$ws.Range("F147").Value = 0
$ws.Range("G147").Value = 0
$ws.Range("H147").Value = 116

# Row 148 now shows Principado de Andorra (its previous, unchanged numbers)
$ws.Range("A148").Value = "Principado de Andorra"
$ws.Range("B148").Value = 3811
$ws.Range("C148").Value = 0
$ws.Range("D148").Value = 2470
$ws.Range("E148").Value = 1278
$ws.Range("F148").Value = 0
$ws.Range("G148").Value = 0
$ws.Range("H148").Value = 63

$wb.Save()
